$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.02
$ws.Range("H2").Value = "2025-03-10 18:29:39"
$ws.Range("I2").Value = "2025-03-10 18:29:39"

# Row 3 updates
$ws.Range("G3").Value = 1.03
$ws.Range("H3").Value = "2025-03-10 18:29:49"
$ws.Range("I3").Value = "2025-03-10 18:29:50"
